$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Long text blocks as single-quoted here-strings (literal, no PS interpolation) ---
$objetivosPt = @'
Dar conhecimentos aos alunos de noções básicas sobre ecologia e impacto das atividades da engenharia sobre o meio ambiente. Conceitos legais e institucionais para o desenvolvimento sustentável.
'@
$programaResumidoPt = @'
Conceitos e Definições. Questões Ambientais. Desenvolvimento Sustentável. Desempenho Ambiental. Processos Ambientais. Norma Ambiental.
'@
$programaPt = @'
CONCEITOS E DEFINIÇÕES. Engenharia Ambiental. Meio Ambiente. Poluição Ambiental. Componentes Ambientais Críticos. QUESTÕES AMBIENTAIS. O Sujeito das Transformações Ambientais. Energia e o Meio Ambiente. Impactos Ambientais nos Três Meios. Equilíbrio Ameaçado. DESENVOLVIMENTO SUSTENTÁVEL. Conceitos Básicos. Aspectos legais. DESEMPENHO AMBIENTAL. Monitoramento Ambiental. Abrangência do Desenvolvimento Sustentável. Definição de Indicadores. Definição de Indicadores Sustentáveis. Indicadores de Desenvolvimento Humano  IDH. Indicadores de Sustentabilidade Ambiental. Controle de Processos Ambientais. PROCESSOS AMBIENTAIS. Controle Processo ETA. Água na Natureza. Caracterização da Água. Indicadores de Qualidade da Água. NORMA AMBIENTAL. Portaria 518. CONAMA 20. Desastre Ecológico.
'@
$bibliografia = @'
1)        BRAGA, B.; HESPANHOL, I.; CONEJO, J. G. L.; MIERZWA, J. C.; BARROS, M. T. L.; SPENCER, M.; PORTO, M.; NUCCI, N.; JULIANO, N.; EIGER, S. Introdução à Engenharia Ambiental: O Desafio do Desenvolvimento Sustentável. Pearson (2ª Edição), 336 p., 2005.2)        VESILIND, P.A.; MORGAN, S. M.; HEINE, L. G. Introdução à Engenharia Ambiental. Cengage (3ª edição), 472 p., 2018.3)        CALIJURI, M. C.; CUNHA, D. G. F. Engenharia Ambiental: Conceitos, Tecnologias e Gestão. Elsevier (1ª Edição), 832 p., 2012.4)        CAPAZ, R. S.; HORTA NOGUEIRA, L. A. Ciências Ambientais para Engenharia. Elsevier (1ª Edição), 252 p., 2014.5)        DAVIS, M. L.; MASTEN, S. J. Princípios de Engenharia Ambiental. Mc Graw Hill Educations (3ª Edição), 872 p., 2016;
'@

# --- 1) Insert two new rows at row 13 (pushes old rows 13-21 down to 15-23) ---
#     New rows inherit formatting/height from the row above (row 12), which matches the
#     target (default height, style 1 on column A) closely enough; we clean up below.
$ws.Rows.Item(13).Resize(2).Insert()

# Remove the stray A13/A14 cells that Insert() may have populated with style only (no value)
$ws.Range("A13").Clear()
$ws.Range("A14").Clear()

# --- 2) Populate the two brand-new rows (13: Docentes responsaveis #1, 14: #2) ---
# Copy formatting from the (already shifted) B15/C15 cells, which still carry the correct
# wrap-text / red-font styles for column B / C, onto the new B13/C13 and B14/C14 cells.
$ws.Range("B15").Copy()
$ws.Range("B13:B14").PasteSpecial(-4122)
$ws.Range("C15").Copy()
$ws.Range("C13:C14").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("B13").Value = "9146830 - Danúbia Caporusso Bargos"
$ws.Range("C13").Value = "9146830 - Danúbia Caporusso Bargos"
$ws.Range("B14").Value = "5464150 - Mariana Consiglio Kasemodel"
$ws.Range("C14").Value = "5464150 - Mariana Consiglio Kasemodel"

# --- 3) Fix the rows whose B/C values were shifted/duplicated in the original file ---
# Row 10: Objetivos (PT)
$ws.Range("B10").Value = $objetivosPt
$ws.Range("C10").Value = $objetivosPt

# Row 15: Programa resumido (PT)
$ws.Range("B15").Value = $programaResumidoPt
$ws.Range("C15").Value = $programaResumidoPt

# Row 17: Programa (PT) -- this row previously had no B/C cells; copy formatting from B18/C18
$ws.Range("B18").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("C18").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("B17").Value = $programaPt
$ws.Range("C17").Value = $programaPt

# Row 20: Método
$ws.Range("B20").Value = "Aulas expositivas com a utilização de recursos de projeções e audiovisual."
$ws.Range("C20").Value = "Aulas expositivas com a utilização de recursos de projeções e audiovisual."

# Row 21: Critério
$ws.Range("B21").Value = "Média ponderada de 2 avaliações escritas com nota final (NF ≥ 5,0)"
$ws.Range("C21").Value = "Média ponderada de 2 avaliações escritas com nota final (NF ≥ 5,0)"

# Row 22: Norma de recuperação
$ws.Range("B22").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova escrita de recuperação a ser aplicada"
$ws.Range("C22").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova escrita de recuperação a ser aplicada"

# Row 23: Bibliografia
$ws.Range("B23").Value = $bibliografia
$ws.Range("C23").Value = $bibliografia

# --- 4) Column layout: split the old shared "min=1 max=2" column group into its own ---
#     entry for column A (same effective width/style as before); mirrors the diff's <cols> change.
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(2).ColumnWidth

Write-Host "Done"
